# Generate Report for Handback
# Adds a new row (row 4) for file "eda6e11f-1f57-46b8-8794-077d6f6e9b52.md"
# to the Overview sheet and the two per-locale detail sheets (zh-cn, de-de),
# mirroring the existing rows for the other two handed-back files.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(4, 1).Value = "eda6e11f-1f57-46b8-8794-077d6f6e9b52.md"
$wsOverview.Cells.Item(4, 2).Value = "e2e\eda6e11f-1f57-46b8-8794-077d6f6e9b52.md"
$wsOverview.Cells.Item(4, 3).Value = ".md"
$wsOverview.Cells.Item(4, 5).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(4, 6).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(4, 7).Value = "2016-08-31 18:51:25"

# Match the formatting used by the existing rows: column B carries the
# "HyperLink" look, column G carries the datetime number format.
$wsOverview.Cells.Item(4, 2).Font.Underline = 2
$wsOverview.Cells.Item(4, 2).Font.Color = 15570276
$wsOverview.Cells.Item(4, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(4, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e1a9e0f7c1a0c6b7f0b8e9a1f2b3c4d5e6f7a8b/e2e/eda6e11f-1f57-46b8-8794-077d6f6e9b52.md", "", "", "e2e\eda6e11f-1f57-46b8-8794-077d6f6e9b52.md")

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Cells.Item(4, 1).Value = "eda6e11f-1f57-46b8-8794-077d6f6e9b52.md"
$wsZhCn.Cells.Item(4, 2).Value = ".md"
$wsZhCn.Cells.Item(4, 3).Value = "Handed back: in sync with en-US"
$wsZhCn.Cells.Item(4, 4).Value = "e2e"
$wsZhCn.Cells.Item(4, 5).Value = "ht"
$wsZhCn.Cells.Item(4, 6).Value = "True"
$wsZhCn.Cells.Item(4, 7).Value = "eda6e11f-1f57-46b8-8794-077d6f6e9b52.b65269c19a2190cc9646a2db85249b5026a960db.zh-cn.xlf"
$wsZhCn.Cells.Item(4, 8).Value = "2016-08-31 18:51:21"
$wsZhCn.Cells.Item(4, 9).Value = "eda6e11f-1f57-46b8-8794-077d6f6e9b52.md"
$wsZhCn.Cells.Item(4, 10).Value = "eda6e11f-1f57-46b8-8794-077d6f6e9b52.b65269c19a2190cc9646a2db85249b5026a960db.zh-cn.xlf"
$wsZhCn.Cells.Item(4, 11).Value = "2016-08-31 18:51:39"
$wsZhCn.Cells.Item(4, 12).Value = ""
$wsZhCn.Cells.Item(4, 13).Value = "True"
$wsZhCn.Cells.Item(4, 14).Value = ""
$wsZhCn.Cells.Item(4, 15).Value = "False"
$wsZhCn.Cells.Item(4, 16).Value = ""

$wsZhCn.Cells.Item(4, 1).Font.Underline = 2
$wsZhCn.Cells.Item(4, 1).Font.Color = 15570276
$wsZhCn.Cells.Item(4, 9).Font.Underline = 2
$wsZhCn.Cells.Item(4, 9).Font.Color = 15570276
$wsZhCn.Cells.Item(4, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(4, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(4, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e1a9e0f7c1a0c6b7f0b8e9a1f2b3c4d5e6f7a8b/e2e/eda6e11f-1f57-46b8-8794-077d6f6e9b52.md", "", "", "eda6e11f-1f57-46b8-8794-077d6f6e9b52.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(4, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9f8e7d6c5b4a3928170615243342516079685940/e2e/eda6e11f-1f57-46b8-8794-077d6f6e9b52.md", "", "", "eda6e11f-1f57-46b8-8794-077d6f6e9b52.md")

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P4"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Cells.Item(4, 1).Value = "eda6e11f-1f57-46b8-8794-077d6f6e9b52.md"
$wsDeDe.Cells.Item(4, 2).Value = ".md"
$wsDeDe.Cells.Item(4, 3).Value = "Handed back: in sync with en-US"
$wsDeDe.Cells.Item(4, 4).Value = "e2e"
$wsDeDe.Cells.Item(4, 5).Value = "ht"
$wsDeDe.Cells.Item(4, 6).Value = "True"
$wsDeDe.Cells.Item(4, 7).Value = "eda6e11f-1f57-46b8-8794-077d6f6e9b52.b65269c19a2190cc9646a2db85249b5026a960db.de-de.xlf"
$wsDeDe.Cells.Item(4, 8).Value = "2016-08-31 18:51:25"
$wsDeDe.Cells.Item(4, 9).Value = "eda6e11f-1f57-46b8-8794-077d6f6e9b52.md"
$wsDeDe.Cells.Item(4, 10).Value = "eda6e11f-1f57-46b8-8794-077d6f6e9b52.b65269c19a2190cc9646a2db85249b5026a960db.de-de.xlf"
$wsDeDe.Cells.Item(4, 11).Value = "2016-08-31 18:51:46"
$wsDeDe.Cells.Item(4, 12).Value = ""
$wsDeDe.Cells.Item(4, 13).Value = "True"
$wsDeDe.Cells.Item(4, 14).Value = ""
$wsDeDe.Cells.Item(4, 15).Value = "False"
$wsDeDe.Cells.Item(4, 16).Value = ""

$wsDeDe.Cells.Item(4, 1).Font.Underline = 2
$wsDeDe.Cells.Item(4, 1).Font.Color = 15570276
$wsDeDe.Cells.Item(4, 9).Font.Underline = 2
$wsDeDe.Cells.Item(4, 9).Font.Color = 15570276
$wsDeDe.Cells.Item(4, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(4, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(4, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e1a9e0f7c1a0c6b7f0b8e9a1f2b3c4d5e6f7a8b/e2e/eda6e11f-1f57-46b8-8794-077d6f6e9b52.md", "", "", "eda6e11f-1f57-46b8-8794-077d6f6e9b52.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(4, 9), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1a2b3c4d5e6f708192a3b4c5d6e7f8091a2b3c4d/e2e/eda6e11f-1f57-46b8-8794-077d6f6e9b52.md", "", "", "eda6e11f-1f57-46b8-8794-077d6f6e9b52.md")

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P4"))

Write-Host "Handback report row added."
